# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-06-25 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-26 Wednesday", 2)

# Update the answers in the table, addressing each cell directly by
# (row, column) to avoid any ambiguity from duplicate/overlapping text
# values between old and new answers.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "646÷9=71, 7"
$t.Cell(1, 2).Range.Text = "744÷3=248, 0"
$t.Cell(1, 3).Range.Text = "877÷5=175, 2"
$t.Cell(1, 4).Range.Text = "936÷8=117, 0"
$t.Cell(1, 5).Range.Text = "730÷5=146, 0"

$t.Cell(5, 1).Range.Text = "683÷6=113, 5"
$t.Cell(5, 2).Range.Text = "654÷3=218, 0"
$t.Cell(5, 3).Range.Text = "931÷4=232, 3"
$t.Cell(5, 4).Range.Text = "409÷2=204, 1"
$t.Cell(5, 5).Range.Text = "348÷8=43, 4"

$t.Cell(9, 1).Range.Text = "470÷6=78, 2"
$t.Cell(9, 2).Range.Text = "471÷8=58, 7"
$t.Cell(9, 3).Range.Text = "872÷9=96, 8"
$t.Cell(9, 4).Range.Text = "651÷8=81, 3"
$t.Cell(9, 5).Range.Text = "951÷3=317, 0"

$t.Cell(13, 1).Range.Text = "616÷6=102, 4"
$t.Cell(13, 2).Range.Text = "274÷5=54, 4"
$t.Cell(13, 3).Range.Text = "709÷8=88, 5"
$t.Cell(13, 4).Range.Text = "983÷5=196, 3"
$t.Cell(13, 5).Range.Text = "722÷4=180, 2"

$t.Cell(17, 1).Range.Text = "845÷6=140, 5"
$t.Cell(17, 2).Range.Text = "232÷4=58, 0"
$t.Cell(17, 3).Range.Text = "665÷4=166, 1"
$t.Cell(17, 4).Range.Text = "165÷3=55, 0"
$t.Cell(17, 5).Range.Text = "195÷9=21, 6"
